# Zeiterfassung_Kasper_Christian.xlsx
# Erfassung der Zeit des Meetings Eintragung
# -> Add a new time-tracking entry row for a meeting ("Vision, Aufgabenteilung")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# The previous row (row 9) used a slightly different date format than the
# other entries (d-mmm instead of m/d/yyyy). Align it with the rest of the
# table before appending the new row, so the new row matches the established
# look of the list.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(9, 1).PasteSpecial(-4122)

# Append the new meeting entry as row 10, reusing the formatting already
# used for the date (column A) and duration (column B) cells above it.
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 2).Copy()
$ws.Cells.Item(10, 2).PasteSpecial(-4122)

$ws.Cells.Item(10, 1).Value2 = 45586
$ws.Cells.Item(10, 2).Value2 = 1.5
$ws.Cells.Item(10, 3).Value = "Besprechung"
$ws.Cells.Item(10, 4).Value = "Vision, Aufgabenteilung"

$ws.Application.CutCopyMode = 0

# Move the active selection to the newly filled-in description cell, as left
# behind by the author after typing the entry.
$ws.Range("D10").Select()
